$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 201

for ($r = 2; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 3).Value2
    $d = [Math]::Round($c * 0.13, 2)
    $e = $c + $d
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}
